$wb = $excel.ActiveWorkbook

# --- Update the Timestamp column (AK2:AK52) on "FBS" - all rows shared the same ---
# --- generated-timestamp string, which was refreshed by the scraper run.       ---
$wsFBS = $wb.Worksheets.Item("FBS")
$newTimestamp = "2025-10-12T18:05:46.827393"
for ($row = 2; $row -le 52; $row++) {
    $wsFBS.Cells.Item($row, 37).Value = $newTimestamp   # column AK = 37
}

# --- Update individual wind_dir_fg (forecast wind direction) cells ---
# FBS sheet (sheet1): Q13 -> SSE, Q25 -> SSW
$wsFBS.Range("Q13").Value = "SSE"
$wsFBS.Range("Q25").Value = "SSW"

# Other sheet (sheet2): S10 -> NNW, S28 -> NNW
$wsOther = $wb.Worksheets.Item("Other")
$wsOther.Range("S10").Value = "NNW"
$wsOther.Range("S28").Value = "NNW"
